$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell values scraped fresh by the GitHub Actions cron job.
# D-column price strings that parse as plain numbers need a leading
# apostrophe so Excel keeps storing them as text (matches the source
# data, which uses literal text cells throughout column D).
$ws.Range("D2").Value = '68.181.40'
$ws.Range("E2").Value = '  +0.72%  '
$ws.Range("D3").Value = '3.835.81'
$ws.Range("E3").Value = '  -0.37%  '
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("D5").Value = '''600.12'
$ws.Range("E5").Value = '  +0.43%  '
$ws.Range("D6").Value = '''171.55'
$ws.Range("E6").Value = '  +3.28%  '
$ws.Range("D7").Value = '3.836.17'
$ws.Range("E7").Value = '  -0.29%  '
$ws.Range("E8").Value = '  +0.01%  '
$ws.Range("E9").Value = '  -0.04%  '
$ws.Range("E10").Value = '  +2.12%  '
$ws.Range("E11").Value = '  +2.58%  '
$ws.Range("E12").Value = '  +1.23%  '
$ws.Range("E13").Value = '  +15.19%  '
$ws.Range("D14").Value = '''36.93'
$ws.Range("E14").Value = '  +0.34%  '
$ws.Range("E15").Value = '  -0.45%  '
$ws.Range("E16").Value = '  -0.29%  '
$ws.Range("D17").Value = '68.257.94'
$ws.Range("E17").Value = '  +0.64%  '
$ws.Range("E18").Value = '  +1.80%  '
$ws.Range("E19").Value = '  +1.53%  '
$ws.Range("E20").Value = '  +0.64%  '
$ws.Range("E21").Value = '  -0.45%  '
$ws.Range("D22").Value = '''468.05'
$ws.Range("E22").Value = '  +1.16%  '
$ws.Range("E24").Value = '  -2.40%  '
$ws.Range("E25").Value = '  +0.40%  '
$ws.Range("D26").Value = '''2.27'
$ws.Range("E26").Value = '  +0.83%  '
$ws.Range("E27").Value = '  +0.14%  '
$ws.Range("D28").Value = '''10.43'
$ws.Range("E28").Value = '  +4.43%  '
$ws.Range("E29").Value = '  +0.08%  '
$ws.Range("E30").Value = '  -0.35%  '
$ws.Range("D31").Value = '3.988.80'
$ws.Range("E31").Value = '  -0.23%  '
$ws.Range("E32").Value = '  +0.09%  '
$ws.Range("E33").Value = '  -0.75%  '
$ws.Range("D34").Value = '''30.99'
$ws.Range("E34").Value = '  -0.05%  '
$ws.Range("D35").Value = '''9.39'
$ws.Range("E35").Value = '  +1.04%  '
$ws.Range("E36").Value = '  -0.73%  '
$ws.Range("D37").Value = '''3.89'
$ws.Range("E37").Value = '  +19.96%  '
$ws.Range("E38").Value = '  +1.18%  '
$ws.Range("D39").Value = '''5.97'
$ws.Range("E39").Value = '  +1.35%  '
$ws.Range("E40").Value = '  -0.03%  '
$ws.Range("E41").Value = '  +0.20%  '
$ws.Range("E43").Value = '  +2.53%  '
$ws.Range("E45").Value = '  +0.71%  '
$ws.Range("D46").Value = '''8.76'
$ws.Range("E46").Value = '  +3.27%  '
$ws.Range("D47").Value = '''417.24'
$ws.Range("E47").Value = '  -2.09%  '
$ws.Range("D48").Value = '''0.000294'
$ws.Range("E48").Value = '  +7.45%  '
$ws.Range("D49").Value = '''46.57'
$ws.Range("E49").Value = '  -1.24%  '
$ws.Range("E50").Value = '  +1.72%  '
$ws.Range("D51").Value = '''141.54'
$ws.Range("E51").Value = '  -1.57%  '
